$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: cells whose text must be preserved as literal text (not
# auto-coerced to numbers) get an explicit "@" text format first.

# Row 8
$ws.Range("C8").Value = 31

# Row 9
$ws.Range("C9").Value = 87
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2"
$ws.Range("E9").Value = "Short point (up to 3 mtr.)"
$ws.Range("F9").Value = 256
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "22272.00"

# Row 10
$ws.Range("C10").Value = 51
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3"
$ws.Range("E10").Value = "Medium point (up to 6 mtr.)"
$ws.Range("F10").Value = 472
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "24072.00"

# Row 11
$ws.Range("A11").Value = "P. point"
$ws.Range("C11").Value = 22
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "4"
$ws.Range("E11").Value = "Long point  (up to 10 mtr.)"
$ws.Range("F11").Value = 662
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "14564.00"

# Row 12
$ws.Range("A12").Value = ""
$ws.Range("C12").Value = 93
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.0"
$ws.Range("E12").Value = "Rewiring of 3/5 pin 6 amp. Light plug point with 1.5 sq. mm nominal size  FR PVC insulated unsheathed flexible copper conductor 1.1 kV grade  and 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper earth conductor 1.1 kV grade(IS:694)   in recessed ISI marked MMS ( IS:9537 P - III ) virgin material  PVC conduit & it's  ISI marked (IS:3419-1988) accessories, 1.2 mm thick  MS box with earth terminal of required size,  6 A  switch, 3/5 pin 6 A socket, 3.0 mm thick ISI marked (IS:2036-1995) phenolic laminated sheet, Al.alloy / Cadmium plated iron/ brass  screws, cup washers, making connections, testing etc. as required.  For specification of copper  Conductor,  Phenolic Laminated sheet's & Electrical/ Wiring accessories refer Chapter E - 04, E - 05 & E - 07 For additional technical parameters of product / work refer Annexure 'A' attached with this BSR"
$ws.Range("F12").Value = 0
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "0.00"

# Row 13
$ws.Range("A13").Value = "P. point"
$ws.Range("C13").Value = 56
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6"
$ws.Range("E13").Value = "On board"
$ws.Range("F13").Value = 136
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "7616.00"

# Row 14
$ws.Range("C14").Value = 86

# Row 15
$ws.Range("C15").Value = 64

# Row 16
$ws.Range("C16").Value = 89

# Row 18 totals
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "68524.00"
$ws.Range("H18").NumberFormat = "@"
$ws.Range("H18").Value = "68524.00"

# Row 20 totals
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "68524.00"
$ws.Range("H20").NumberFormat = "@"
$ws.Range("H20").Value = "68524.00"

